$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.337.58"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "2.605.32"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D5").Value = "588.45"
$ws.Range("E5").Value = "  +6.33%  "
$ws.Range("D6").Value = "143.01"
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").Value = "2.610.28"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("E13").Value = "  +4.49%  "
$ws.Range("D14").Value = "3.067.79"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "24.66"
$ws.Range("E15").Value = "  +6.38%  "
$ws.Range("D16").Value = "60.340.91"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").Value = "2.611.52"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "11.38"
$ws.Range("E19").Value = "  +9.77%  "
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").Value = "347.60"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").Value = "  +6.10%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +9.13%  "
$ws.Range("D25").Value = "63.08"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +7.76%  "
$ws.Range("D29").Value = "0.0₃0794"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  +10.95%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "6.37"
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Value = "163.55"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("D34").Value = "19.52"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "0.982"
$ws.Range("D37").Value = "1.24"
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("D38").Value = "1.62"
$ws.Range("E38").Value = "  +10.06%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +6.33%  "
$ws.Range("D41").Value = "310.49"
$ws.Range("E41").Value = "  +7.51%  "
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "135.71"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "0.0994"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "19.77"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "5.01"
$ws.Range("E47").Value = "  +10.54%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.604"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").Value = "20.17"
$ws.Range("E50").Value = "  +8.22%  "
$ws.Range("D51").Value = "0.0241"
$ws.Range("E51").Value = "  +3.12%  "
